# Daily User Impact Status - append the next day's row (row 22) and
# move the active selection to it, matching the author's upload.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New day: 2025-11-05 (serial 45966) with its metric row.
$ws.Range("A22").Value = 45966
$ws.Range("B22").Value = 5600
$ws.Range("C22").Value = 4370
$ws.Range("D22").Value = 4046
$ws.Range("E22").Value = 258
$ws.Range("F22").Value = 44
$ws.Range("G22").Value = 17
$ws.Range("H22").Value = 5
$ws.Range("I22").Value = 0

# Move the selected/active cell to the newly added row, as in the saved file.
$excel.Goto($ws.Range("A22:I22"))
